$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.770.42"
$ws.Range("E2").Value = "  +0.92%  "
$ws.Range("D3").Value = "'3.273.33"
$ws.Range("E3").Value = "  -1.84%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "'576.88"
$ws.Range("E5").Value = "  -1.34%  "
$ws.Range("D6").Value = "'171.89"
$ws.Range("E6").Value = "  -7.19%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "'3.265.48"
$ws.Range("E9").Value = "  -1.88%  "
$ws.Range("E10").Value = "  -5.32%  "
$ws.Range("D11").Value = "'0.568"
$ws.Range("E11").Value = "  -2.24%  "
$ws.Range("D12").Value = "'44.84"
$ws.Range("E12").Value = "  -4.68%  "
$ws.Range("D13").Value = "'0.0000265"
$ws.Range("E13").Value = "  -1.39%  "
$ws.Range("D14").Value = "'690.63"
$ws.Range("E14").Value = "  +2.82%  "
$ws.Range("D15").Value = "'3.803.19"
$ws.Range("E15").Value = "  -1.68%  "
$ws.Range("E16").Value = "  -3.58%  "
$ws.Range("D17").Value = "'66.878.84"
$ws.Range("E17").Value = "  +0.82%  "
$ws.Range("E18").Value = "  +0.73%  "
$ws.Range("D19").Value = "'3.274.05"
$ws.Range("E19").Value = "  -1.77%  "
$ws.Range("D20").Value = "'17.14"
$ws.Range("E20").Value = "  -4.09%  "
$ws.Range("D21").Value = "'10.65"
$ws.Range("E21").Value = "  -3.92%  "
$ws.Range("D22").Value = "'0.881"
$ws.Range("E22").Value = "  -1.87%  "
$ws.Range("D23").Value = "'16.85"
$ws.Range("E23").Value = "  -4.99%  "
$ws.Range("D24").Value = "'5.19"
$ws.Range("E24").Value = "  +2.99%  "
$ws.Range("D25").Value = "'99.02"
$ws.Range("E25").Value = "  -3.62%  "
$ws.Range("E26").Value = "  -4.00%  "
$ws.Range("E27").Value = "  -5.38%  "
$ws.Range("E28").Value = "  +2.96%  "
$ws.Range("E29").Value = "  -4.50%  "
$ws.Range("E30").Value = "  -2.76%  "
$ws.Range("D31").Value = "'6.58"
$ws.Range("E31").Value = "  -4.49%  "
$ws.Range("D32").Value = "'577.68"
$ws.Range("E32").Value = "  -5.39%  "
$ws.Range("D33").Value = "'10.80"
$ws.Range("E33").Value = "  -2.63%  "
$ws.Range("D34").Value = "'3.816.35"
$ws.Range("E34").Value = "  -1.05%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("E36").Value = "  -3.51%  "
$ws.Range("D37").Value = "'55.02"
$ws.Range("E37").Value = "  -1.78%  "
$ws.Range("D38").Value = "'3.28"
$ws.Range("E38").Value = "  -15.86%  "
$ws.Range("D39").Value = "'0.128"
$ws.Range("E39").Value = "  -0.33%  "
$ws.Range("D40").Value = "'3.38"
$ws.Range("E40").Value = "  -0.94%  "
$ws.Range("E41").Value = "  -4.71%  "
$ws.Range("D42").Value = "'31.27"
$ws.Range("E42").Value = "  -4.76%  "
$ws.Range("D43").Value = "'0.0₃0658"
$ws.Range("E43").Value = "  -6.26%  "
$ws.Range("E44").Value = "  -4.30%  "
$ws.Range("D45").Value = "'2.95"
$ws.Range("E45").Value = "  -7.42%  "
$ws.Range("E46").Value = "  -4.28%  "
$ws.Range("E47").Value = "  -0.13%  "
$ws.Range("E48").Value = "  -1.49%  "
$ws.Range("D49").Value = "'2.53"
$ws.Range("E49").Value = "  -0.92%  "
$ws.Range("E50").Value = "  +3.38%  "
$ws.Range("D51").Value = "'129.01"
$ws.Range("E51").Value = "  -0.50%  "
